# Auto-generated update to profit/price figures across the Leve sheets.
# Each block updates the currentAveragePrice / LevePrice / LeveProfit columns
# (H:N) for one row, keyed by its Leve Item ID (column G), matching a refreshed
# market-data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Item ID 27772)
$ws.Range("H28").Value = 1190.75
$ws.Range("I28").Value = 1203.7142
$ws.Range("J28").Value = 1100
$ws.Range("K28").Value = 1203.7142
$ws.Range("L28").Value = 1100
$ws.Range("M28").Value = -718.7141999999999
$ws.Range("N28").Value = -2070

# Row 40 (Item ID 5505)
$ws.Range("H40").Value = 3174.2593
$ws.Range("I40").Value = 2153.2942
$ws.Range("K40").Value = 2153.2942
$ws.Range("M40").Value = -1978.2942

# Row 86 (Item ID 12603)
$ws.Range("H86").Value = 8273.532999999999
$ws.Range("I86").Value = 8152.1
$ws.Range("J86").Value = 8516.4
$ws.Range("K86").Value = 8152.1
$ws.Range("L86").Value = 8516.4
$ws.Range("M86").Value = -7029.1
$ws.Range("N86").Value = -10762.4

# Row 89 (Item ID 12603)
$ws.Range("H89").Value = 8273.532999999999
$ws.Range("I89").Value = 8152.1
$ws.Range("J89").Value = 8516.4
$ws.Range("K89").Value = 40760.5
$ws.Range("L89").Value = 42582
$ws.Range("M89").Value = -35144.5
$ws.Range("N89").Value = -53814

# Row 98 (Item ID 36237)
$ws.Range("H98").Value = 1953.3334
$ws.Range("I98").Value = 1953.3334
$ws.Range("K98").Value = 1953.3334
$ws.Range("M98").Value = -455.3334

# Row 99 (Item ID 19883)
$ws.Range("H99").Value = 333.33334
$ws.Range("I99").Value = 400
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 600
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -3596

# Row 101 (Item ID 19884)
$ws.Range("H101").Value = 3675.25
$ws.Range("I101").Value = 1567
$ws.Range("K101").Value = 4701
$ws.Range("M101").Value = -3079

# Row 107 (Item ID 27766)
$ws.Range("H107").Value = 2042.25
$ws.Range("I107").Value = 1364.2727
$ws.Range("K107").Value = 1364.2727
$ws.Range("M107").Value = 555.7273

# Row 115 (Item ID 27957)
$ws.Range("H115").Value = 1149.4286
$ws.Range("I115").Value = 609.2
$ws.Range("K115").Value = 1827.6
$ws.Range("M115").Value = -260.6000000000001

# Row 118 (Item ID 27958)
$ws.Range("H118").Value = 785.4545000000001
$ws.Range("I118").Value = 836.5
$ws.Range("J118").Value = 275
$ws.Range("K118").Value = 2509.5
$ws.Range("L118").Value = 825
$ws.Range("M118").Value = -852.5
$ws.Range("N118").Value = -4139

# Row 121 (Item ID 39731)
$ws.Range("H121").Value = 854
$ws.Range("J121").Value = 854
$ws.Range("L121").Value = 2562
$ws.Range("N121").Value = -6056

# Row 122 (Item ID 36237)
$ws.Range("H122").Value = 1953.3334
$ws.Range("I122").Value = 1953.3334
$ws.Range("K122").Value = 5860.0002
$ws.Range("M122").Value = -3410.0002

# Row 127 (Item ID 36114)
$ws.Range("H127").Value = 288
$ws.Range("I127").Value = 288
$ws.Range("K127").Value = 864
$ws.Range("M127").Value = 4096

# Row 129 (Item ID 36115)
$ws.Range("H129").Value = 1929.1111
$ws.Range("I129").Value = 1819
$ws.Range("J129").Value = 2017.2
$ws.Range("K129").Value = 5457
$ws.Range("L129").Value = 6051.6
$ws.Range("M129").Value = -457
$ws.Range("N129").Value = -16051.6

# Row 132 (Item ID 44049)
$ws.Range("H132").Value = 3174.2222
$ws.Range("I132").Value = 3174.2222
$ws.Range("K132").Value = 9522.6666
$ws.Range("M132").Value = -6992.6666

# Row 137 (Item ID 44013)
$ws.Range("H137").Value = 2764
$ws.Range("I137").Value = 1198.5
$ws.Range("J137").Value = 3333.2727
$ws.Range("K137").Value = 3595.5
$ws.Range("L137").Value = 9999.8181
$ws.Range("M137").Value = -1045.5
$ws.Range("N137").Value = -15099.8181

# Row 138 (Item ID 44169)
$ws.Range("H138").Value = 3809.3865
$ws.Range("J138").Value = 3779.7778
$ws.Range("L138").Value = 11339.3334
$ws.Range("N138").Value = -21619.3334

$ws = $wb.Worksheets.Item("ARM")
# Row 3 (Item ID 2494)
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

# Row 32 (Item ID 44147)
$ws.Range("H32").Value = 8153.8237
$ws.Range("I32").Value = 4783.364
$ws.Range("K32").Value = 4783.364
$ws.Range("M32").Value = -4496.364

# Row 132 (Item ID 43997)
$ws.Range("H132").Value = 2501582.2
$ws.Range("I132").Value = 2779188.8
$ws.Range("J132").Value = 3124.5
$ws.Range("K132").Value = 8337566.399999999
$ws.Range("L132").Value = 9373.5
$ws.Range("M132").Value = -8335036.399999999
$ws.Range("N132").Value = -14433.5

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Item ID 5092)
$ws.Range("H22").Value = 232
$ws.Range("J22").Value = 400
$ws.Range("L22").Value = 400
$ws.Range("N22").Value = -746

# Row 86 (Item ID 12526)
$ws.Range("H86").Value = 4402
$ws.Range("I86").Value = 4250.75
$ws.Range("K86").Value = 4250.75
$ws.Range("M86").Value = -3127.75

# Row 89 (Item ID 12526)
$ws.Range("H89").Value = 4402
$ws.Range("I89").Value = 4250.75
$ws.Range("K89").Value = 21253.75
$ws.Range("M89").Value = -15637.75

# Row 134 (Item ID 43998)
$ws.Range("H134").Value = 62501890
$ws.Range("I134").Value = 62501890
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 187505670
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -187503135
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Item ID 44023)
$ws.Range("H31").Value = 4430
$ws.Range("I31").Value = 5928.7144
$ws.Range("J31").Value = 1432.5714
$ws.Range("K31").Value = 5928.7144
$ws.Range("L31").Value = 1432.5714
$ws.Range("M31").Value = -5633.7144
$ws.Range("N31").Value = -2022.5714

# Row 34 (Item ID 44023)
$ws.Range("H34").Value = 4430
$ws.Range("I34").Value = 5928.7144
$ws.Range("J34").Value = 1432.5714
$ws.Range("K34").Value = 5928.7144
$ws.Range("L34").Value = 1432.5714
$ws.Range("M34").Value = -5726.7144
$ws.Range("N34").Value = -1836.5714

# Row 58 (Item ID 44021)
$ws.Range("H58").Value = 21747466
$ws.Range("I58").Value = 26324964
$ws.Range("K58").Value = 26324964
$ws.Range("M58").Value = -26324761

# Row 62 (Item ID 12580)
$ws.Range("H62").Value = 2826
$ws.Range("I62").Value = 2268.3333
$ws.Range("J62").Value = 4499
$ws.Range("K62").Value = 2268.3333
$ws.Range("L62").Value = 4499
$ws.Range("M62").Value = -1644.3333
$ws.Range("N62").Value = -5747

# Row 65 (Item ID 12580)
$ws.Range("H65").Value = 2826
$ws.Range("I65").Value = 2268.3333
$ws.Range("J65").Value = 4499
$ws.Range("K65").Value = 11341.6665
$ws.Range("L65").Value = 22495
$ws.Range("M65").Value = -8221.666499999999
$ws.Range("N65").Value = -28735

# Row 94 (Item ID 32934)
$ws.Range("H94").Value = 2144.5
$ws.Range("I94").Value = 2423.3333
$ws.Range("J94").Value = 1865.6666
$ws.Range("K94").Value = 2423.3333
$ws.Range("L94").Value = 1865.6666
$ws.Range("M94").Value = -1972.3333
$ws.Range("N94").Value = -2767.6666

# Row 136 (Item ID 44021)
$ws.Range("H136").Value = 21747466
$ws.Range("I136").Value = 26324964
$ws.Range("K136").Value = 78974892
$ws.Range("M136").Value = -78972342

# Row 141 (Item ID 43345)
$ws.Range("H141").Value = 86712
$ws.Range("J141").Value = 84304.664
$ws.Range("L141").Value = 84304.664
$ws.Range("N141").Value = -94664.664

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (Item ID 4847)
$ws.Range("H2").Value = 627.4286
$ws.Range("I2").Value = 97
$ws.Range("J2").Value = 715.8333
$ws.Range("K2").Value = 582
$ws.Range("L2").Value = 4294.9998
$ws.Range("M2").Value = -469
$ws.Range("N2").Value = -4520.9998

# Row 38 (Item ID 4860)
$ws.Range("H38").Value = 184.06667
$ws.Range("I38").Value = 181.77777
$ws.Range("K38").Value = 545.33331
$ws.Range("M38").Value = -198.33331

# Row 40 (Item ID 4827)
$ws.Range("H40").Value = 5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 20
$ws.Range("N40").Value = -158
$ws.Range("M40").ClearContents()

# Row 137 (Item ID 44088)
$ws.Range("H137").Value = 2102.889
$ws.Range("I137").Value = 1885.4
$ws.Range("J137").Value = 2374.75
$ws.Range("K137").Value = 5656.200000000001
$ws.Range("L137").Value = 7124.25
$ws.Range("M137").Value = -556.2000000000007
$ws.Range("N137").Value = -17324.25

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Item ID 12521)
$ws.Range("H80").Value = 2249.5
$ws.Range("I80").Value = 2249.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2249.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1251.5
$ws.Range("N80").ClearContents()

# Row 83 (Item ID 12521)
$ws.Range("H83").Value = 2249.5
$ws.Range("I83").Value = 2249.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 11247.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -6255.5
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 93 (Item ID 19993)
$ws.Range("H93").Value = 2952.375
$ws.Range("I93").Value = 2302.7144
$ws.Range("K93").Value = 2302.7144
$ws.Range("M93").Value = -1054.7144

# Row 136 (Item ID 44060)
$ws.Range("H136").Value = 4659
$ws.Range("J136").Value = 6660.6665
$ws.Range("L136").Value = 19981.9995
$ws.Range("N136").Value = -25081.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Item ID 12596)
$ws.Range("H81").Value = 1424.25
$ws.Range("I81").Value = 1424.25
$ws.Range("K81").Value = 2848.5
$ws.Range("M81").Value = -1787.5

# Row 84 (Item ID 12596)
$ws.Range("H84").Value = 1424.25
$ws.Range("I84").Value = 1424.25
$ws.Range("K84").Value = 14242.5
$ws.Range("M84").Value = -8938.5

# Row 122 (Item ID 36208)
$ws.Range("H122").Value = 1336.6666
$ws.Range("I122").Value = 1336.6666
$ws.Range("K122").Value = 4009.9998
$ws.Range("M122").Value = -1559.9998

# Row 132 (Item ID 44029)
$ws.Range("H132").Value = 125001730
$ws.Range("I132").Value = 125001730
$ws.Range("K132").Value = 375005190
$ws.Range("M132").Value = -375002660

# Row 136 (Item ID 44031)
$ws.Range("H136").Value = 20834428
$ws.Range("I136").Value = 22728262
$ws.Range("K136").Value = 68184786
$ws.Range("M136").Value = -68182236
